$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 19533.396
$ws.Range("I132").Value = 2934.4524
$ws.Range("K132").Value = 8803.3572
$ws.Range("M132").Value = -6273.3572

$ws.Range("H135").Value = 15153027
$ws.Range("I135").Value = 900.2222
$ws.Range("J135").Value = 83337600
$ws.Range("K135").Value = 8101.999800000001
$ws.Range("L135").Value = 750038400
$ws.Range("M135").Value = -5566.999800000001
$ws.Range("N135").Value = -750043470

$ws.Range("H137").Value = 3618.4883
$ws.Range("I137").Value = 933.0741
$ws.Range("K137").Value = 2799.2223
$ws.Range("M137").Value = -249.2223000000004

$ws.Range("H138").Value = 1397.32
$ws.Range("I138").Value = 636.48834
$ws.Range("J138").Value = 1971.2808
$ws.Range("K138").Value = 1909.46502
$ws.Range("L138").Value = 5913.8424
$ws.Range("M138").Value = 3230.53498
$ws.Range("N138").Value = -16193.8424

$ws.Range("H141").Value = 3438.0938
$ws.Range("I141").Value = 1269.7
$ws.Range("J141").Value = 7052.0835
$ws.Range("K141").Value = 3809.1
$ws.Range("L141").Value = 21156.2505
$ws.Range("M141").Value = 1370.9
$ws.Range("N141").Value = -31516.2505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H61").Value = 1102.0465
$ws.Range("I61").Value = 836.40625
$ws.Range("K61").Value = 836.40625
$ws.Range("M61").Value = -624.40625

$ws.Range("H74").Value = 1455.1731
$ws.Range("I74").Value = 1229.1702
$ws.Range("K74").Value = 1229.1702
$ws.Range("M74").Value = -355.1702

$ws.Range("H77").Value = 1455.1731
$ws.Range("I77").Value = 1229.1702
$ws.Range("K77").Value = 6145.851000000001
$ws.Range("M77").Value = -1777.851000000001

$ws.Range("H122").Value = 1911.4615
$ws.Range("I122").Value = 1923.1818
$ws.Range("K122").Value = 5769.5454
$ws.Range("M122").Value = -3319.5454

$ws.Range("H132").Value = 2332.3872
$ws.Range("I132").Value = 1328.7894
$ws.Range("J132").Value = 3921.4167
$ws.Range("K132").Value = 3986.3682
$ws.Range("L132").Value = 11764.2501
$ws.Range("M132").Value = -1456.3682
$ws.Range("N132").Value = -16824.2501

$ws.Range("H136").Value = 1102.0465
$ws.Range("I136").Value = 836.40625
$ws.Range("K136").Value = 2509.21875
$ws.Range("M136").Value = 40.78125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3413.0356
$ws.Range("I134").Value = 3156
$ws.Range("K134").Value = 9468
$ws.Range("M134").Value = -6933

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4142.684
$ws.Range("I31").Value = 1578.88
$ws.Range("J31").Value = 5058.3286
$ws.Range("K31").Value = 1578.88
$ws.Range("L31").Value = 5058.3286
$ws.Range("M31").Value = -1283.88
$ws.Range("N31").Value = -5648.3286

$ws.Range("H32").Value = 936.3333
$ws.Range("I32").Value = 936.3333
$ws.Range("K32").Value = 936.3333
$ws.Range("M32").Value = -620.3333

$ws.Range("H34").Value = 4142.684
$ws.Range("I34").Value = 1578.88
$ws.Range("J34").Value = 5058.3286
$ws.Range("K34").Value = 1578.88
$ws.Range("L34").Value = 5058.3286
$ws.Range("M34").Value = -1376.88
$ws.Range("N34").Value = -5462.3286

$ws.Range("H58").Value = 1786.1613
$ws.Range("I58").Value = 1535.2727
$ws.Range("K58").Value = 1535.2727
$ws.Range("M58").Value = -1332.2727

$ws.Range("H99").Value = 2646.3
$ws.Range("I99").Value = 1941.8
$ws.Range("J99").Value = 3350.8
$ws.Range("K99").Value = 1941.8
$ws.Range("L99").Value = 3350.8
$ws.Range("M99").Value = -443.8
$ws.Range("N99").Value = -6346.8

$ws.Range("H122").Value = 100899.836
$ws.Range("I122").Value = 240799.6
$ws.Range("K122").Value = 722398.8
$ws.Range("M122").Value = -719948.8

$ws.Range("H126").Value = 2646.3
$ws.Range("I126").Value = 1941.8
$ws.Range("J126").Value = 3350.8
$ws.Range("K126").Value = 5825.4
$ws.Range("L126").Value = 10052.4
$ws.Range("M126").Value = -3355.4
$ws.Range("N126").Value = -14992.4

$ws.Range("H132").Value = 59318.44
$ws.Range("I132").Value = 2002.375
$ws.Range("K132").Value = 6007.125
$ws.Range("M132").Value = -3477.125

$ws.Range("H134").Value = 401480.62
$ws.Range("I134").Value = 1065.8966
$ws.Range("K134").Value = 3197.6898
$ws.Range("M134").Value = -662.6898000000001

$ws.Range("H136").Value = 1786.1613
$ws.Range("I136").Value = 1535.2727
$ws.Range("K136").Value = 4605.8181
$ws.Range("M136").Value = -2055.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 10869653
$ws.Range("I2").Value = 42.357143
$ws.Range("J2").Value = 27777936
$ws.Range("K2").Value = 254.142858
$ws.Range("L2").Value = 166667616
$ws.Range("M2").Value = -141.142858
$ws.Range("N2").Value = -166667842

$ws.Range("H4").Value = 890.1429000000001
$ws.Range("I4").Value = 85.8
$ws.Range("J4").Value = 2901
$ws.Range("K4").Value = 257.4
$ws.Range("L4").Value = 8703
$ws.Range("M4").Value = -145.4
$ws.Range("N4").Value = -8927

$ws.Range("H8").Value = 2061
$ws.Range("I8").Value = 2061
$ws.Range("K8").Value = 6183
$ws.Range("M8").Value = -6044

$ws.Range("H9").Value = 100008030
$ws.Range("J9").Value = 8793.5
$ws.Range("L9").Value = 26380.5
$ws.Range("N9").Value = -26828.5

$ws.Range("H14").Value = 29.4
$ws.Range("I14").Value = 29.4
$ws.Range("K14").Value = 88.19999999999999
$ws.Range("M14").Value = 84.80000000000001

$ws.Range("H20").Value = 2340
$ws.Range("I20").Value = 425
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 1275
$ws.Range("L20").Value = 30000
$ws.Range("M20").Value = -1048
$ws.Range("N20").Value = -30454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H21").Value = 7100
$ws.Range("J21").Value = 7100
$ws.Range("L21").Value = 7100
$ws.Range("N21").Value = -7446

$ws.Range("H30").Value = 7100
$ws.Range("J30").Value = 7100
$ws.Range("L30").Value = 7100
$ws.Range("N30").Value = -7310

$ws.Range("H132").Value = 3106.9312
$ws.Range("I132").Value = 2217.3333
$ws.Range("J132").Value = 4562.636
$ws.Range("K132").Value = 6651.999899999999
$ws.Range("L132").Value = 13687.908
$ws.Range("M132").Value = -4121.999899999999
$ws.Range("N132").Value = -18747.908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3380.7
$ws.Range("I7").Value = 2600.3333
$ws.Range("J7").Value = 4551.25
$ws.Range("K7").Value = 2600.3333
$ws.Range("L7").Value = 4551.25
$ws.Range("M7").Value = -2488.3333
$ws.Range("N7").Value = -4775.25

$ws.Range("H22").Value = 779.93335
$ws.Range("I22").Value = 714.1429000000001
$ws.Range("J22").Value = 837.5
$ws.Range("K22").Value = 714.1429000000001
$ws.Range("L22").Value = 837.5
$ws.Range("M22").Value = -419.1429000000001
$ws.Range("N22").Value = -1427.5

$ws.Range("H27").Value = 779.93335
$ws.Range("I27").Value = 714.1429000000001
$ws.Range("J27").Value = 837.5
$ws.Range("K27").Value = 714.1429000000001
$ws.Range("L27").Value = 837.5
$ws.Range("M27").Value = -607.1429000000001
$ws.Range("N27").Value = -1051.5

$ws.Range("H34").Value = 13500
$ws.Range("J34").Value = 13500
$ws.Range("L34").Value = 13500
$ws.Range("N34").Value = -13844

$ws.Range("H40").Value = 3335.652
$ws.Range("I40").Value = 2198
$ws.Range("J40").Value = 5468.75
$ws.Range("K40").Value = 2198
$ws.Range("L40").Value = 5468.75
$ws.Range("M40").Value = -2062
$ws.Range("N40").Value = -5740.75

$ws.Range("H61").Value = 1942.5333
$ws.Range("I61").Value = 2198.5833
$ws.Range("J61").Value = 918.3333
$ws.Range("K61").Value = 2198.5833
$ws.Range("L61").Value = 918.3333
$ws.Range("M61").Value = -1996.5833
$ws.Range("N61").Value = -1322.3333

$ws.Range("H113").Value = 1942.5333
$ws.Range("I113").Value = 2198.5833
$ws.Range("J113").Value = 918.3333
$ws.Range("K113").Value = 2198.5833
$ws.Range("L113").Value = 918.3333
$ws.Range("M113").Value = -28.58329999999978
$ws.Range("N113").Value = -5258.3333

$ws.Range("H126").Value = 3380.7
$ws.Range("I126").Value = 2600.3333
$ws.Range("J126").Value = 4551.25
$ws.Range("K126").Value = 7800.999899999999
$ws.Range("L126").Value = 13653.75
$ws.Range("M126").Value = -5330.999899999999
$ws.Range("N126").Value = -18593.75

$ws.Range("H132").Value = 1949.4546
$ws.Range("I132").Value = 1422.3934
$ws.Range("K132").Value = 4267.1802
$ws.Range("M132").Value = -1737.1802

$ws.Range("H136").Value = 2030.0869
$ws.Range("I136").Value = 1676.1177
$ws.Range("K136").Value = 5028.3531
$ws.Range("M136").Value = -2478.3531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1464.8541
$ws.Range("I132").Value = 1153.1765
$ws.Range("K132").Value = 3459.5295
$ws.Range("M132").Value = -929.5295000000001

$ws.Range("H136").Value = 294922.72
$ws.Range("I136").Value = 345454.3
$ws.Range("K136").Value = 1036362.9
$ws.Range("M136").Value = -1033812.9
